# Update the COVID-19 Valais daily data sheet with newly published figures
# for the dates 2021-09-15 through 2021-09-20 (rows 569-573), plus a small
# correction to row 564 and fill-ins for the previously-empty rows 570-573.
# Columns: A=Date, B=Cumul cas positifs (formula), C=Nb nouveaux cas positifs,
# D=Nb nouvelles admissions, E=Patients SI, F=Patients intubes,
# G=Patients hospitalises hors SI, H=Total hospitalisations (formula),
# I=Nb nouvelles sorties, J=Cumul deces (formula), K=Nb nouveaux deces (formula),
# L=Nb nouveaux deces hopital, M=Nb nouveaux deces extra-hospitaliers.
#
# Columns L and M are formatted as Text ("@") in this sheet. Writing a plain
# number into a Text-formatted cell stores it as text, so for those two
# columns we briefly switch the cell to General, write the number, then
# restore the original Text format (this keeps the cell's style/format the
# same as before while still storing a real numeric value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-NumericDeathCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "General"
    $cell.Value = $value
    $cell.NumberFormat = "@"
}

# Row 564: one additional extra-hospital death recorded
Set-NumericDeathCell 564 13 1

# Row 568: revised new-cases and intubated counts
$ws.Cells.Item(568, 3).Value = 87
$ws.Cells.Item(568, 6).Value = 5

# Row 569: revised new-cases count; one hospital death recorded
$ws.Cells.Item(569, 3).Value = 72
Set-NumericDeathCell 569 12 1

# Row 570: newly published data
$ws.Cells.Item(570, 3).Value = 53
$ws.Cells.Item(570, 5).Value = 9
$ws.Cells.Item(570, 6).Value = 2
$ws.Cells.Item(570, 7).Value = 10
Set-NumericDeathCell 570 12 0
Set-NumericDeathCell 570 13 1

# Row 571: newly published data
$ws.Cells.Item(571, 3).Value = 18
$ws.Cells.Item(571, 5).Value = 8
$ws.Cells.Item(571, 6).Value = 2
$ws.Cells.Item(571, 7).Value = 9
Set-NumericDeathCell 571 12 0
Set-NumericDeathCell 571 13 0

# Row 572: newly published data
$ws.Cells.Item(572, 3).Value = 13
$ws.Cells.Item(572, 5).Value = 8
$ws.Cells.Item(572, 6).Value = 2
$ws.Cells.Item(572, 7).Value = 12
Set-NumericDeathCell 572 12 0
Set-NumericDeathCell 572 13 0

# Row 573: newly published data
$ws.Cells.Item(573, 3).Value = 1
$ws.Cells.Item(573, 5).Value = 8
$ws.Cells.Item(573, 6).Value = 2
$ws.Cells.Item(573, 7).Value = 13
Set-NumericDeathCell 573 12 0
Set-NumericDeathCell 573 13 0

# Reset the frozen-pane scroll position back to the top of the data and
# select the header row, as the author did before re-uploading the file.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("A1:M1").Select()
